$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect D2:D51 from Excel auto-converting numeric-looking text (e.g. "20.30") into real numbers
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.207.14"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "1.632.92"
$ws.Range("E3").Value = "  -0.73%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "216.69"
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("D6").Value = "0.516"
$ws.Range("E6").Value = "  +1.04%  "
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("D8").Value = "0.255"
$ws.Range("E8").Value = "  -0.44%  "
$ws.Range("E9").Value = "  -0.84%  "
$ws.Range("D10").Value = "20.30"
$ws.Range("E10").Value = "  +1.49%  "
$ws.Range("D11").Value = "0.0849"
$ws.Range("E11").Value = "  +0.57%  "
$ws.Range("D12").Value = "1.618.11"
$ws.Range("E12").Value = "  -1.48%  "
$ws.Range("E13").Value = "  -0.04%  "
$ws.Range("D14").Value = "0.544"
$ws.Range("E14").Value = "  +1.17%  "
$ws.Range("D15").Value = "27.204.38"
$ws.Range("E15").Value = "  +0.49%  "
$ws.Range("D16").Value = "64.84"
$ws.Range("E16").Value = "  -3.75%  "
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("D18").Value = "215.90"
$ws.Range("E18").Value = "  -1.34%  "
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("D20").Value = "6.91"
$ws.Range("E20").Value = "  +0.79%  "
$ws.Range("E21").Value = "  -0.61%  "
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("D23").Value = "9.08"
$ws.Range("E23").Value = "  -1.23%  "
$ws.Range("D24").Value = "148.48"
$ws.Range("E24").Value = "  +0.41%  "
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("D26").Value = "7.30"
$ws.Range("E26").Value = "  -1.26%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").Value = "15.58"
$ws.Range("E28").Value = "  -1.12%  "
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("E30").Value = "  -0.80%  "
$ws.Range("D31").Value = "3.39"
$ws.Range("E31").Value = "  +0.54%  "
$ws.Range("D32").Value = "3.02"
$ws.Range("E32").Value = "  -0.70%  "
$ws.Range("D33").Value = "1.317.07"
$ws.Range("E33").Value = "  +4.09%  "
$ws.Range("E34").Value = "  -1.60%  "
$ws.Range("E35").Value = "  -0.38%  "
$ws.Range("D36").Value = "0.0175"
$ws.Range("E36").Value = "  -1.30%  "
$ws.Range("D37").Value = "0.849"
$ws.Range("E37").Value = "  +0.88%  "
$ws.Range("D38").Value = "0.539"
$ws.Range("E38").Value = "  -0.42%  "
$ws.Range("E39").Value = "  -0.31%  "
$ws.Range("D40").Value = "2.27"
$ws.Range("E40").Value = "  +2.11%  "
$ws.Range("D41").Value = "0.805"
$ws.Range("E41").Value = "  -0.57%  "
$ws.Range("D42").Value = "63.78"
$ws.Range("E42").Value = "  +2.45%  "
$ws.Range("D43").Value = "1.771.20"
$ws.Range("E43").Value = "  -0.82%  "
$ws.Range("E44").Value = "  -2.48%  "
$ws.Range("D45").Value = "90.89"
$ws.Range("E45").Value = "  -1.12%  "
$ws.Range("E46").Value = "  +0.43%  "
$ws.Range("E47").Value = "  +5.47%  "
$ws.Range("D48").Value = "0.798"
$ws.Range("E48").Value = "  +18.06%  "
$ws.Range("D49").Value = "0.0516"
$ws.Range("E49").Value = "  +0.79%  "
$ws.Range("D50").Value = "7.57"
$ws.Range("E50").Value = "  -1.43%  "
$ws.Range("E51").Value = "  -0.47%  "

# Restore the default (Normal) style so cells keep their original, unstyled appearance
$ws.Range("D2:D51").Style = "Normal"
